# Fix a duplicated-dot typo in the ParserException handler code sample:
#   "errorHandler..reportError" -> "errorHandler.reportError"
# The original markup stored this as three separate runs
# ("errorHandler" / ".." / "reportError"); the fix merges them into a
# single run containing the corrected text while preserving the
# surrounding runs (leading whitespace run and trailing "(e);" run)
# untouched.

$p = $ppt.ActivePresentation

$oldText = "errorHandler..reportError"
$newText = "errorHandler.reportError"

$targetSlide = $null
$targetShape = $null

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $txt = $shp.TextFrame.TextRange.Text
                if ($txt.IndexOf($oldText) -ge 0) {
                    $targetSlide = $sl
                    $targetShape = $shp
                }
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange
$fullText = $tr.Text
$charIndex = $fullText.IndexOf($oldText)

$startPos = $charIndex + 1
$length = $oldText.Length

$run = $tr.Characters($startPos, $length)
$run.Text = $newText
